# db/dummydata/st_local_charges.xlsx -- "for obi and dejan"
#
# 1. Header D1 "SERVICE LEVEL" -> "SERVICE_LEVEL" (underscore, matching the
#    other ALL_CAPS_WITH_UNDERSCORES headers).
# 2. Column K (RATE_BASIS) had a handful of rows using inconsistent-case
#    duplicates of the standard rate-basis codes ("per_CBM", "per_shipment",
#    "per_Shipment") instead of the canonical "PER_CBM" / "PER_SHIPMENT"
#    used everywhere else. Normalize those cells to the canonical values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Header rename --------------------------------------------------------
$ws.Range("D1").Value = "SERVICE_LEVEL"

# -- RATE_BASIS normalization (column K) ----------------------------------
# "per_CBM" -> "PER_CBM"
$perCbmRows = 31,33,35,37,40,42,46,48,49,58
foreach ($r in $perCbmRows) {
    $ws.Cells.Item($r, 11).Value = "PER_CBM"
}

# "per_shipment" / "per_Shipment" -> "PER_SHIPMENT"
$perShipmentRows = 39,41,43,44
foreach ($r in $perShipmentRows) {
    $ws.Cells.Item($r, 11).Value = "PER_SHIPMENT"
}
